$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the header row (row 1: "Name", "Age") so the data shifts up.
$ws.Rows.Item(1).Delete()
